$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Content changes ("documentatie einde van de week" updates)
# Order matters for shared-string append order in the saved OOXML.
$ws.Range("C50").Value = "Bug fix (probleem met kruispunt lichtwaardes en B = control) + belbin test in teamcontract"
$ws.Range("C47").Value = "Documentatie tot nu toe op hu intranet zetten"
$ws.Range("C41").Value = "activity diagrams"

# Update the view's active selection to match the saved workbook state
$ws.Range("F53").Select()
